$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-looking strings in column D stay as plain text (not auto-converted
# to Excel date serials) across the whole range we are about to populate.
$ws.Range("D2:D8").NumberFormat = "@"

# New column L ("OCR usado") - clone the bold/centered header style from A1,
# then overwrite with the new header text.
$ws.Range("A1").Copy($ws.Range("L1"))
$ws.Range("L1").Value = "OCR usado"

# Row data: Resolución, Nombre de norma, Descripción, Fecha, Archivo,
# publication_type_id, category_id, Nombre de Archivo, Descripción del
# documento, RUTA TEMP (file name), OCR usado (bool)
$rows = @(
    @{ Row=2; A="431-2025-MPH"; File="documento_ejemplo copy 2.pdf" },
    @{ Row=3; A="431-2025-MPH"; File="documento_ejemplo copy 3.pdf" },
    @{ Row=4; A="431-2025-MPH"; File="documento_ejemplo copy 4.pdf" },
    @{ Row=5; A="431-2025-MPH"; File="documento_ejemplo copy 5.pdf" },
    @{ Row=6; A="431-2025-MPH"; File="documento_ejemplo copy 6.pdf" },
    @{ Row=7; A="431-2025-MPH"; File="documento_ejemplo copy.pdf" },
    @{ Row=8; A="431-2025-MPH"; File="documento_ejemplo.pdf" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = ";"
    $ws.Cells.Item($row, 3).Value = ";"
    $ws.Cells.Item($row, 4).Value = "05/09/2025"
    $ws.Cells.Item($row, 5).Value = "https://tustorage.municipalidad.gob.pe/archivos/" + $r.File
    $ws.Cells.Item($row, 6).Value = 159
    $ws.Cells.Item($row, 7).Value = 54
    $ws.Cells.Item($row, 8).Value = "RESOLUCION 431-2025-MPH"
    $ws.Cells.Item($row, 10).Value = "Documento 431-2025-MPH"
    $ws.Cells.Item($row, 11).Value = $r.File
    $ws.Cells.Item($row, 12).Value = $true
}
